# Weekly update: insert a new Berenjena price record for
# "Vega Monumental Concepción" as row 67 in the data table, pushing the
# existing rows 67-77 down to 68-78 (dimension grows from R77 to R78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 67:77 down one row to make room for the new record.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new weekly sample.
$ws.Cells.Item(67, 1).Value = 11
$ws.Cells.Item(67, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(67, 3).Value = "Bíobío"
$ws.Cells.Item(67, 4).Value = 44637
$ws.Cells.Item(67, 5).Value = 8
$ws.Cells.Item(67, 6).Value = 100112001
$ws.Cells.Item(67, 7).Value = "Berenjena"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 180
$ws.Cells.Item(67, 11).Value = 9000
$ws.Cells.Item(67, 12).Value = 10000
$ws.Cells.Item(67, 13).Value = 9556
$ws.Cells.Item(67, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(67, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(67, 16).Value = 159
$ws.Cells.Item(67, 17).Value = 60
$ws.Cells.Item(67, 18).Value = "Hortaliza"
